$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1082.0834
$ws.Range("J41").Value = 524.5714
$ws.Range("L41").Value = 524.5714
$ws.Range("N41").Value = -1404.5714
$ws.Range("H51").Value = 2879.6
$ws.Range("I51").Value = 2625
$ws.Range("J51").Value = 3114.6155
$ws.Range("K51").Value = 2625
$ws.Range("L51").Value = 3114.6155
$ws.Range("M51").Value = -2141
$ws.Range("N51").Value = -4082.6155
$ws.Range("H58").Value = 3733.9092
$ws.Range("J58").Value = 4200
$ws.Range("L58").Value = 12600
$ws.Range("N58").Value = -12900
$ws.Range("H112").Value = 33305.855
$ws.Range("J112").Value = 35251.97
$ws.Range("L112").Value = 105755.91
$ws.Range("N112").Value = -107971.91
$ws.Range("H132").Value = 3084.2
$ws.Range("I132").Value = 3090.2144
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 9270.643199999999
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -6740.643199999999
$ws.Range("N132").Value = -14060
$ws.Range("H138").Value = 6581985
$ws.Range("J138").Value = 7695585
$ws.Range("L138").Value = 23086755
$ws.Range("N138").Value = -23097035
$ws.Range("H141").Value = 2797.5
$ws.Range("I141").Value = 2797.5
$ws.Range("K141").Value = 8392.5
$ws.Range("M141").Value = -3212.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6789.5776
$ws.Range("I32").Value = 3206.1025
$ws.Range("K32").Value = 3206.1025
$ws.Range("M32").Value = -2919.1025
$ws.Range("H61").Value = 4363.4253
$ws.Range("I61").Value = 3128.375
$ws.Range("K61").Value = 3128.375
$ws.Range("M61").Value = -2916.375
$ws.Range("H74").Value = 52592.863
$ws.Range("I74").Value = 72691.14
$ws.Range("J74").Value = 17420.875
$ws.Range("K74").Value = 72691.14
$ws.Range("L74").Value = 17420.875
$ws.Range("M74").Value = -71817.14
$ws.Range("N74").Value = -19168.875
$ws.Range("H76").Value = 100000
$ws.Range("J76").Value = 100000
$ws.Range("L76").Value = 100000
$ws.Range("N76").Value = -100676
$ws.Range("H77").Value = 52592.863
$ws.Range("I77").Value = 72691.14
$ws.Range("J77").Value = 17420.875
$ws.Range("K77").Value = 363455.7
$ws.Range("L77").Value = 87104.375
$ws.Range("M77").Value = -359087.7
$ws.Range("N77").Value = -95840.375
$ws.Range("H79").Value = 100000
$ws.Range("J79").Value = 100000
$ws.Range("L79").Value = 100000
$ws.Range("N79").Value = -102340
$ws.Range("H110").Value = 4428.0464
$ws.Range("I110").Value = 4327.9443
$ws.Range("K110").Value = 4327.9443
$ws.Range("M110").Value = -2282.9443
$ws.Range("H132").Value = 3407.54
$ws.Range("I132").Value = 3219.366
$ws.Range("K132").Value = 9658.098
$ws.Range("M132").Value = -7128.098
$ws.Range("H136").Value = 4363.4253
$ws.Range("I136").Value = 3128.375
$ws.Range("K136").Value = 9385.125
$ws.Range("M136").Value = -6835.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 9675.75
$ws.Range("I86").Value = 9629.571
$ws.Range("J86").Value = 9999
$ws.Range("K86").Value = 9629.571
$ws.Range("L86").Value = 9999
$ws.Range("M86").Value = -8506.571
$ws.Range("N86").Value = -12245
$ws.Range("H89").Value = 9675.75
$ws.Range("I89").Value = 9629.571
$ws.Range("J89").Value = 9999
$ws.Range("K89").Value = 48147.855
$ws.Range("L89").Value = 49995
$ws.Range("M89").Value = -42531.855
$ws.Range("N89").Value = -61227
$ws.Range("H107").Value = 1437.3077
$ws.Range("I107").Value = 1299.091
$ws.Range("K107").Value = 1299.091
$ws.Range("M107").Value = 620.9090000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37707.793
$ws.Range("I31").Value = 50451.24
$ws.Range("K31").Value = 50451.24
$ws.Range("M31").Value = -50156.24
$ws.Range("H34").Value = 37707.793
$ws.Range("I34").Value = 50451.24
$ws.Range("K34").Value = 50451.24
$ws.Range("M34").Value = -50249.24
$ws.Range("H132").Value = 4233.7334
$ws.Range("I132").Value = 3530.1538
$ws.Range("K132").Value = 10590.4614
$ws.Range("M132").Value = -8060.4614
$ws.Range("H134").Value = 19993.072
$ws.Range("I134").Value = 12812.3
$ws.Range("J134").Value = 37945
$ws.Range("K134").Value = 38436.89999999999
$ws.Range("L134").Value = 113835
$ws.Range("M134").Value = -35901.89999999999
$ws.Range("N134").Value = -118905

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 5861246.5
$ws.Range("I4").Value = 6041292.5
$ws.Range("K4").Value = 18123877.5
$ws.Range("M4").Value = -18123765.5
$ws.Range("H18").Value = 923.25
$ws.Range("I18").Value = 886.6667
$ws.Range("K18").Value = 2660.0001
$ws.Range("M18").Value = -2491.0001
$ws.Range("H38").Value = 83.85714
$ws.Range("I38").Value = 71.666664
$ws.Range("J38").Value = 93
$ws.Range("K38").Value = 214.999992
$ws.Range("L38").Value = 279
$ws.Range("M38").Value = 132.000008
$ws.Range("N38").Value = -973
$ws.Range("H50").Value = 2491.5
$ws.Range("I50").Value = 974.5
$ws.Range("J50").Value = 3250
$ws.Range("K50").Value = 2923.5
$ws.Range("L50").Value = 9750
$ws.Range("M50").Value = -2442.5
$ws.Range("N50").Value = -10712
$ws.Range("H53").Value = 2491.5
$ws.Range("I53").Value = 974.5
$ws.Range("J53").Value = 3250
$ws.Range("K53").Value = 2923.5
$ws.Range("L53").Value = 9750
$ws.Range("M53").Value = -2442.5
$ws.Range("N53").Value = -10712
$ws.Range("H141").Value = 125848.375
$ws.Range("I141").Value = 964.8570999999999
$ws.Range("J141").Value = 1000033
$ws.Range("K141").Value = 2894.5713
$ws.Range("L141").Value = 3000099
$ws.Range("M141").Value = 2285.4287
$ws.Range("N141").Value = -3010459

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 38000.25
$ws.Range("I7").Value = 2002
$ws.Range("J7").Value = 49999.668
$ws.Range("K7").Value = 2002
$ws.Range("L7").Value = 49999.668
$ws.Range("N7").Value = -50223.668
$ws.Range("M7").Value = -1890
$ws.Range("H8").Value = 38000.25
$ws.Range("I8").Value = 2002
$ws.Range("J8").Value = 49999.668
$ws.Range("K8").Value = 2002
$ws.Range("L8").Value = 49999.668
$ws.Range("N8").Value = -50277.668
$ws.Range("M8").Value = -1863
$ws.Range("H132").Value = 3516
$ws.Range("I132").Value = 2931.762
$ws.Range("J132").Value = 4631.364
$ws.Range("K132").Value = 8795.286
$ws.Range("L132").Value = 13894.092
$ws.Range("M132").Value = -6265.286
$ws.Range("N132").Value = -18954.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2002.76
$ws.Range("I61").Value = 1863.6522
$ws.Range("J61").Value = 3602.5
$ws.Range("K61").Value = 1863.6522
$ws.Range("L61").Value = 3602.5
$ws.Range("M61").Value = -1661.6522
$ws.Range("N61").Value = -4006.5
$ws.Range("H113").Value = 2002.76
$ws.Range("I113").Value = 1863.6522
$ws.Range("J113").Value = 3602.5
$ws.Range("K113").Value = 1863.6522
$ws.Range("L113").Value = 3602.5
$ws.Range("M113").Value = 306.3478
$ws.Range("N113").Value = -7942.5
$ws.Range("H132").Value = 3380.4375
$ws.Range("I132").Value = 3234.652
$ws.Range("K132").Value = 9703.956
$ws.Range("M132").Value = -7173.956
$ws.Range("H136").Value = 3809.0588
$ws.Range("I136").Value = 3528.2222
$ws.Range("J136").Value = 4125
$ws.Range("K136").Value = 10584.6666
$ws.Range("L136").Value = 12375
$ws.Range("M136").Value = -8034.6666
$ws.Range("N136").Value = -17475

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7081.423
$ws.Range("I81").Value = 9776.462
$ws.Range("J81").Value = 4386.385
$ws.Range("K81").Value = 19552.924
$ws.Range("L81").Value = 8772.77
$ws.Range("M81").Value = -18491.924
$ws.Range("N81").Value = -10894.77
$ws.Range("H84").Value = 7081.423
$ws.Range("I84").Value = 9776.462
$ws.Range("J84").Value = 4386.385
$ws.Range("K84").Value = 97764.62
$ws.Range("L84").Value = 43863.85000000001
$ws.Range("M84").Value = -92460.62
$ws.Range("N84").Value = -54471.85000000001
$ws.Range("H113").Value = 674.86365
$ws.Range("I113").Value = 718.375
$ws.Range("K113").Value = 2155.125
$ws.Range("M113").Value = 14.875
$ws.Range("H123").Value = 99999
$ws.Range("J123").Value = 99999
$ws.Range("L123").Value = 99999
$ws.Range("N123").Value = -109799
$ws.Range("H136").Value = 1948.7037
$ws.Range("I136").Value = 1302.4117
$ws.Range("K136").Value = 3907.2351
$ws.Range("M136").Value = -1357.2351

